$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before D (shifts D:K -> E:L), carrying formatting along
$ws.Columns("D").Insert(-4161)

# Copy formatting (number formats/styles) from the now-shifted E column into
# the new, blank D column for each of the three data blocks in the sheet.
$ws.Range("E7:E35").Copy()
$ws.Range("D7:D35").PasteSpecial(-4122)
$ws.Range("E38:E77").Copy()
$ws.Range("D38:D77").PasteSpecial(-4122)
$ws.Range("E80:E102").Copy()
$ws.Range("D80:D102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Populate the new column D with the latest fiscal-year figures
$ws.Range("D7").Value = 43465
$ws.Range("D8").Value = 859600
$ws.Range("D9").Value = 390500
$ws.Range("D10").Value = 469100
$ws.Range("D12").Value = 38000
$ws.Range("D13").Value = 0
$ws.Range("D14").Value = 9000
$ws.Range("D15").Value = 23200
$ws.Range("D17").Value = 788300
$ws.Range("D18").Value = 71300
$ws.Range("D20").Value = 0
$ws.Range("D21").Value = 133100
$ws.Range("D22").Value = 20700
$ws.Range("D23").Value = 50700
$ws.Range("D24").Value = 10700
$ws.Range("D25").Value = 0
$ws.Range("D26").Value = 40000
$ws.Range("D27").Value = 40000
$ws.Range("D28").Value = 0
$ws.Range("D29").Value = 900
$ws.Range("D30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("D32").Value = 0
$ws.Range("D33").Value = 40900
$ws.Range("D34").Value = 0
$ws.Range("D35").Value = 40900
$ws.Range("D38").Value = 43465
$ws.Range("D41").Value = 17500
$ws.Range("D42").Value = 0
$ws.Range("D43").Value = 181600
$ws.Range("D44").Value = 154600
$ws.Range("D45").Value = 20700
$ws.Range("D46").Value = 374400
$ws.Range("D47").Value = 0
$ws.Range("D48").Value = 113200
$ws.Range("D49").Value = 813600
$ws.Range("D50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("D52").Value = 67900
$ws.Range("D53").Value = 0
$ws.Range("D54").Value = 1369100
$ws.Range("D57").Value = 53500
$ws.Range("D58").Value = 18300
$ws.Range("D59").Value = 89100
$ws.Range("D60").Value = 160900
$ws.Range("D61").Value = 438600
$ws.Range("D62").Value = 107400
$ws.Range("D63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("D66").Value = 706900
$ws.Range("D68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("D72").Value = 464900
$ws.Range("D73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("D76").Value = 662300
$ws.Range("D77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("D81").Value = 40900
$ws.Range("D83").Value = 61800
$ws.Range("D84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("D89").Value = 74700
$ws.Range("D91").Value = -16500
$ws.Range("D92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("D94").Value = -16500
$ws.Range("D96").Value = -22400
$ws.Range("D97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("D100").Value = -72300
$ws.Range("D101").Value = -1000
$ws.Range("D102").Value = -15100

Write-Host "Inserted new column D and populated latest-year figures"
